# Applies the scheduled-runner market-price refresh to the per-job
# profit tables on every class sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# Each sheet has the same 14-column layout:
#   H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#   K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ
#
# Some rows don't carry every column (M/N are omitted when not applicable),
# so a handful of updates add or clear a cell rather than just overwrite its
# value - that's reproduced exactly below using ClearContents()/Value=.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        $ws,
        [int]$Row,
        [hashtable]$Values   # column letter -> new value, or $null to clear the cell
    )
    foreach ($col in $Values.Keys) {
        $cellRef = "$col$Row"
        $val = $Values[$col]
        if ($null -eq $val) {
            $ws.Range($cellRef).ClearContents()
        } else {
            $ws.Range($cellRef).Value = $val
        }
    }
}

# ---------------------------------------------------------------- ALC ----
$ws = $wb.Worksheets.Item("ALC")
Set-Row $ws 16  @{ H=3000;      I=0;        J=3000;      K=0;          L=3000;      M=$null;  N=-3460 }
Set-Row $ws 39  @{ H=434.25806; I=254.41667; J=547.8421;  K=763.25001;  L=1643.5263; M=-467.25001; N=-2235.5263 }
Set-Row $ws 62  @{ H=5681.1;    I=2268.3333;              K=2268.3333;               M=-1644.3333 }
Set-Row $ws 65  @{ H=5681.1;    I=2268.3333;              K=11341.6665;              M=-8221.6665 }
Set-Row $ws 103 @{ H=1167.3;    I=964.6667; J=1254.1428; K=2894.0001;  L=3762.4284; M=-2308.0001; N=-4934.428400000001 }
Set-Row $ws 137 @{ H=5828.375;              J=12239.833;               L=36719.499;               N=-41819.499 }
Set-Row $ws 141 @{ H=2895.1333; I=2901.9285; J=2800;      K=8705.7855; L=8400;      M=-3525.7855; N=-18760 }

# ---------------------------------------------------------------- ARM ----
$ws = $wb.Worksheets.Item("ARM")
Set-Row $ws 32  @{ H=10872109;  I=11906483; J=11178;     K=11906483;  L=11178;     M=-11906196; N=-11752 }
Set-Row $ws 56  @{ H=3000;      I=3000;     J=0;         K=3000;      L=0;         M=-2258;  N=$null }
Set-Row $ws 102 @{ H=23866.46;  I=23866.46;              K=23866.46;               M=-22244.46 }
Set-Row $ws 110 @{ H=1732;      I=1732;                  K=1732;                   M=313 }

# ---------------------------------------------------------------- BSM ----
$ws = $wb.Worksheets.Item("BSM")
Set-Row $ws 99  @{ H=9466.071;  I=10227.083;             K=10227.083;              M=-8729.083 }
Set-Row $ws 107 @{ H=3821.5;    I=2394;                  K=2394;                   M=-474 }

# ---------------------------------------------------------------- CRP ----
$ws = $wb.Worksheets.Item("CRP")
Set-Row $ws 31  @{ H=353136.78; I=4187.2764;             K=4187.2764;              M=-3892.2764 }
Set-Row $ws 34  @{ H=353136.78; I=4187.2764;             K=4187.2764;              M=-3985.2764 }
Set-Row $ws 58  @{ H=6006;      I=6006;                  K=6006;                   M=-5803 }
Set-Row $ws 62  @{ H=4499.8335; I=3999.6667;             K=3999.6667;              M=-3375.6667 }
Set-Row $ws 65  @{ H=4499.8335; I=3999.6667;             K=19998.3335;             M=-16878.3335 }
Set-Row $ws 132 @{ H=1960;      I=1657.2727; J=3625;     K=4971.8181; L=10875;     M=-2441.8181; N=-15935 }
Set-Row $ws 136 @{ H=6006;      I=6006;                  K=18018;                  M=-15468 }

# ---------------------------------------------------------------- CUL ----
$ws = $wb.Worksheets.Item("CUL")
Set-Row $ws 50  @{ H=477.77777;             J=477.77777;               L=1433.33331;              N=-2395.33331 }
Set-Row $ws 53  @{ H=477.77777;             J=477.77777;               L=1433.33331;              N=-2395.33331 }
Set-Row $ws 74  @{ H=14333.333; I=0;        J=14333.333; K=0;         L=42999.999; M=$null;  N=-45121.999 }
Set-Row $ws 77  @{ H=14333.333; I=0;        J=14333.333; K=0;         L=128999.997; M=$null; N=-139607.997 }
Set-Row $ws 86  @{ H=431.6;     I=200;                   K=600;                    M=586 }
Set-Row $ws 89  @{ H=431.6;     I=200;                   K=1800;                   M=4128 }
Set-Row $ws 105 @{ H=0;                     J=0;                       L=0;                       N=$null }

# ---------------------------------------------------------------- GSM ----
$ws = $wb.Worksheets.Item("GSM")
Set-Row $ws 59  @{ H=7119.25;               J=7119.25;                 L=7119.25;                 N=-8285.25 }
Set-Row $ws 102 @{ H=4573.864;  I=3164.0625;             K=3164.0625;              M=-1542.0625 }
Set-Row $ws 132 @{ H=76925464;  I=76925464;              K=230776392;              M=-230773862 }
Set-Row $ws 134 @{ H=60000;                 J=60000;                   L=180000;                  N=-185070 }

# ---------------------------------------------------------------- LTW ----
$ws = $wb.Worksheets.Item("LTW")
Set-Row $ws 122 @{ H=6359;      I=5398.636; J=9000;      K=16195.908; L=27000;     M=-13745.908; N=-31900 }
Set-Row $ws 132 @{ H=119939.65; I=64935.562;             K=194806.686;             M=-192276.686 }

# ---------------------------------------------------------------- WVR ----
$ws = $wb.Worksheets.Item("WVR")
Set-Row $ws 2   @{ H=78500;                 J=50000;                   L=50000;                   N=-50224 }
Set-Row $ws 5   @{ H=4102919.5;             J=5003649.5;               L=5003649.5;               N=-5003873.5 }
Set-Row $ws 132 @{ H=2183.8;    I=2004.3182;             K=6012.9546;              M=-3482.9546 }

Write-Host "Behemoth_Profits: market price refresh applied."
